$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "838 0670 4556"
$ws.Range("B3").Value = "863 8387 5679"
$ws.Range("C2").Value = 778530
$ws.Range("C3").Value = 329714
$ws.Range("E3").Value = "doan training"
$ws.Range("E2").Value = "practice partners"
$ws.Range("A2").Value = 0.71527777777777779

$ws.Range("A3").Select() | Out-Null
